# Populate the daily driver report with a header row + one data row, and
# style the header row: bold font, thin box border, centered/top aligned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$headers = @("name", "employee_id", "asset", "arrival", "status", "division", "job_title")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# ---- Data row (row 2) ----
$record = @("Roger Doddy", "DODROG", "PT-07S", "04:45 AM", "On Time", "TEXDIST", "Select Maintenance Employee")
for ($col = 1; $col -le $record.Length; $col++) {
    $ws.Cells.Item(2, $col).Value = $record[$col - 1]
}

# ---- Header styling: bold, centered/top aligned, thin box border ----
# Build the combined style on A1 first (a single cell keeps the style table
# from growing an orphaned xf per property write), then fan it out to the
# rest of the header row via a format-only copy/paste.
$headerCell = $ws.Range("A1")
$headerCell.Font.Bold = $true
$headerCell.Borders.LineStyle = 1
$headerCell.HorizontalAlignment = -4108
$headerCell.VerticalAlignment = -4160

$headerCell.Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
